$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the new "Comments"-style values in column K ---
# Cells are written in an order chosen so that the workbook's shared-string
# table grows in the same sequence as the target file.
$ws.Range("K39").Value = "PED being used for access"
$ws.Range("K47").Value = "PED being used for access"
$ws.Range("K111").Value = "PED being used for access"
$ws.Range("K113").Value = "PED being used for access"
$ws.Range("K168").Value = "PSD opening"
$ws.Range("K173").Value = "PSD opening"
$ws.Range("K8").Value = "General fault alarm as a result of loss of power"
$ws.Range("K10").Value = "Door forced open"
$ws.Range("K17").Value = "Door forced open"
$ws.Range("K114").Value = "Door forced open"
$ws.Range("K146").Value = "Door forced open"
$ws.Range("K147").Value = "Door forced open"
$ws.Range("K48").Value = "Isolation of lift causes this alarm"
$ws.Range("K49").Value = "Isolation of lift causes this alarm"
$ws.Range("K54").Value = "Isolation of lift causes this alarm"
$ws.Range("K60").Value = "Isolation of lift causes this alarm"
$ws.Range("K149").Value = "PSD isolated during a blockade"
$ws.Range("K169").Value = "PSD isolated during a blockade"
$ws.Range("K170").Value = "PSD isolated during a blockade"
$ws.Range("K171").Value = "PSD isolated during a blockade"
$ws.Range("K172").Value = "PSD isolated during a blockade"
$ws.Range("K185").Value = "PSD isolated during a blockade"
$ws.Range("K186").Value = "PSD isolated during a blockade"
$ws.Range("K187").Value = "PSD isolated during a blockade"
$ws.Range("K61").Value = "Isolation of panel"
$ws.Range("K258").Value = "Not currently connected to Sunshine"
$ws.Range("K259").Value = "Not currently connected to Sunshine"
$ws.Range("K260").Value = "Not currently connected to Sunshine"
$ws.Range("K261").Value = "Not currently connected to Sunshine"
$ws.Range("K262").Value = "Not currently connected to Sunshine"
$ws.Range("K277").Value = "Not currently connected to Sunshine"
$ws.Range("K279").Value = "Not currently connected to Sunshine"
$ws.Range("K188").Value = "Real alarm to be looked at"
$ws.Range("K9").Value = "CYP holding doors open on site"
$ws.Range("K18").Value = "CYP holding doors open on site"
$ws.Range("K38").Value = "CYP holding doors open on site"
$ws.Range("K115").Value = "CYP holding doors open on site"
$ws.Range("K118").Value = "CYP holding doors open on site"
$ws.Range("K148").Value = "CYP holding doors open on site"
$ws.Range("K189").Value = "CYP holding doors open on site"
$ws.Range("K190").Value = "CYP holding doors open on site"

# --- Column K width ---
$ws.Columns("K").ColumnWidth = 40.2

# --- AutoFilter: extend range to include column K, drop stale sort state ---
$ws.AutoFilterMode = $false
$ws.Range("A1:K280").AutoFilter(10, @("ARN"), 7)

# --- Defined name _FilterDatabase must track the same new range ---
$wb.Names.Item("Sheet1!_FilterDatabase").RefersTo = "=Sheet1!`$A`$1:`$K`$280"

# --- Selection / active cell moves to K49 ---
$ws.Range("K49").Select()
